$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 63
$ws.Cells.Item($row, 1).Value = "RandomForest"
$ws.Cells.Item($row, 2).Value = "{'n_estimators': 150, 'max_depth': 20, 'random_state': 42}"
$ws.Cells.Item($row, 3).Value = 0.07107472561565627
$ws.Cells.Item($row, 4).Value = 1
$ws.Cells.Item($row, 5).Value = 0.227799580037436
$ws.Cells.Item($row, 6).Value = 0.9469226425748165
